# Add "start_date" / "end_date" columns to the institution table (Tableau1),
# between "phone" and "tag_ids", and populate a handful of sample values -
# mirrors the author's "added : institution valididy start and end date" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$tbl = $ws.ListObjects.Item(1)

# Remember the width used by the "phone" column (F) so the two freshly
# inserted columns (G:H) end up the same width as their neighbour instead of
# falling back to the sheet default width.
$phoneWidth = $ws.Range("F1").EntireColumn.ColumnWidth

# Detach the existing table definition - we'll rebuild it below once the new
# header cells are in place so the column list/order matches the sheet.
$tbl.Unlist()

# Shift "tag_ids"/"doc_ids" (and everything else in G:H) two columns to the
# right, opening up G:H for the new fields.
$ws.Range("G1:H1").EntireColumn.Insert()
$ws.Range("G1:H1").EntireColumn.ColumnWidth = $phoneWidth

# New header cells.
$ws.Range("G1").Value = "start_date"
$ws.Range("H1").Value = "end_date"

# A few sample values, matching what the source workbook now has.
$ws.Range("G5").Value = 2012
$ws.Range("H5").Value = 2023
$ws.Range("G8").Value = "2010/10"
$ws.Range("H10").Value = "2021/04"

# Re-create the table over the expanded range so column names/order are
# picked up straight from the header row.
$newTbl = $ws.ListObjects.Add(1, $ws.Range("A1:J70"), 0, 1)
$newTbl.Name = "Tableau1"
$newTbl.TableStyle = "TableStyleMedium9"

# Restore the selection to roughly where the author ended up.
$ws.Range("H11").Select()
